$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Sending cluster"="ECs" rows (2-4) are removed from the dataset; the
# remaining "Sending cluster"="MuSCs" rows move up to rows 2-4 and their
# derived-specificity figures are recalculated against the new (smaller)
# dataset (ECs no longer contributes).

$ws.Range("2:4").EntireRow.Delete()

# Row 2 (was row 5): MuSCs -> Clec4g/Lag3 -> ECs
$ws.Range("G2").Value = 0.011591
$ws.Range("H2").Value = 0.034773
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 7.574702666666667
$ws.Range("N2").Value = 22.724108
$ws.Range("O2").Value = 0.2360813295275979
$ws.Range("P2").Value = 0.2360813295275979
$ws.Range("Q2").Value = 0.08779837860933333
$ws.Range("R2").Value = 0.790185407484
$ws.Range("S2").Value = 0.2360813295275979
$ws.Range("T2").Value = 0.2360813295275979

# Row 3 (was row 6): MuSCs -> Clec4g/Lag3 -> FAPs
$ws.Range("G3").Value = 0.011591
$ws.Range("H3").Value = 0.034773
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.5879438355171306
$ws.Range("P3").Value = 0.5879438355171307
$ws.Range("Q3").Value = 0.2186556453873333
$ws.Range("R3").Value = 1.967900808486
$ws.Range("S3").Value = 0.5879438355171306
$ws.Range("T3").Value = 0.5879438355171307

# Row 4 (was row 7): MuSCs -> Clec4g/Lag3 -> MuSCs
$ws.Range("G4").Value = 0.011591
$ws.Range("H4").Value = 0.034773
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("O4").Value = 0.1759748349552714
$ws.Range("P4").Value = 0.1759748349552714
$ws.Range("Q4").Value = 0.06544484147066666
$ws.Range("R4").Value = 0.5890035732359999
$ws.Range("S4").Value = 0.1759748349552714
$ws.Range("T4").Value = 0.1759748349552714
